$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.119.56'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '3.517.41'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '''570.96'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('D6').Value = '''182.62'
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('D7').Value = '3.509.60'
$ws.Range('E7').Value = '  -1.22%  '
$ws.Range('D8').Value = '''0.614'
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').Value = '''0.186'
$ws.Range('E10').Value = '  +5.53%  '
$ws.Range('E11').Value = '  -2.92%  '
$ws.Range('D12').Value = '''53.77'
$ws.Range('E12').Value = '  -3.59%  '
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '''9.46'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = '4.062.18'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = '''19.23'
$ws.Range('E16').Value = '  -3.10%  '
$ws.Range('D17').Value = '3.510.74'
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').Value = '69.106.81'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D21').Value = '''539.21'
$ws.Range('E21').Value = '  +13.99%  '
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('D23').Value = '''20.65'
$ws.Range('E23').Value = '  +6.73%  '
$ws.Range('D24').Value = '''5.00'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('D26').Value = '''94.33'
$ws.Range('E26').Value = '  +6.98%  '
$ws.Range('D27').Value = '''10.98'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D28').Value = '''2.90'
$ws.Range('E28').Value = '  -4.30%  '
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('D30').Value = '''31.53'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('D31').Value = '''7.28'
$ws.Range('E31').Value = '  -4.40%  '
$ws.Range('D32').Value = '''12.69'
$ws.Range('E32').Value = '  +5.45%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.114'
$ws.Range('E33').Value = '  -4.26%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '''63.88'
$ws.Range('E34').Value = '  -2.63%  '
$ws.Range('D35').Value = '''568.05'
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('D36').Value = '''38.18'
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('D37').Value = '''3.08'
$ws.Range('E37').Value = '  +8.22%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').Value = '''0.398'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').Value = '0.0₃0762'
$ws.Range('E40').Value = '  -4.21%  '
$ws.Range('D41').Value = '''3.14'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('E42').Value = '  -4.55%  '
$ws.Range('E43').Value = '  -4.17%  '
$ws.Range('D44').Value = '''3.52'
$ws.Range('E44').Value = '  +6.05%  '
$ws.Range('D45').Value = '''2.96'
$ws.Range('E45').Value = '  -4.33%  '
$ws.Range('D46').Value = '3.188.16'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').Value = '''9.17'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('E49').Value = '  -2.10%  '
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('D51').Value = '''136.44'
$ws.Range('E51').Value = '  -0.62%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D51').Style = "Normal"
